$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 516.678
$ws.Range("J17").Value = 529.7143
$ws.Range("L17").Value = 1589.1429
$ws.Range("N17").Value = -1925.1429

$ws.Range("H18").Value = 1001
$ws.Range("I18").Value = 1001
$ws.Range("K18").Value = 1001
$ws.Range("M18").Value = -717

$ws.Range("H41").Value = 1094.9375
$ws.Range("I41").Value = 944.5
$ws.Range("J41").Value = 1245.375
$ws.Range("K41").Value = 944.5
$ws.Range("L41").Value = 1245.375
$ws.Range("M41").Value = -504.5
$ws.Range("N41").Value = -2125.375

$ws.Range("H105").Value = 114000
$ws.Range("J105").Value = 114000
$ws.Range("L105").Value = 114000
$ws.Range("N105").Value = -120988

$ws.Range("H111").Value = 3206.3635
$ws.Range("I111").Value = 3127.1
$ws.Range("K111").Value = 9381.299999999999
$ws.Range("M111").Value = -6314.299999999999

$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = 254

$ws.Range("H117").Value = 59880
$ws.Range("J117").Value = 59880
$ws.Range("L117").Value = 59880
$ws.Range("N117").Value = -69058

$ws.Range("H120").Value = 53170
$ws.Range("J120").Value = 53170
$ws.Range("L120").Value = 53170
$ws.Range("N120").Value = -62846

$ws.Range("H137").Value = 3623.361
$ws.Range("I137").Value = 1775.091
$ws.Range("J137").Value = 6527.7856
$ws.Range("K137").Value = 5325.272999999999
$ws.Range("L137").Value = 19583.3568
$ws.Range("M137").Value = -2775.272999999999
$ws.Range("N137").Value = -24683.3568


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1012.4
$ws.Range("I2").Value = 1028.2858
$ws.Range("K2").Value = 1028.2858
$ws.Range("M2").Value = -915.2858000000001

$ws.Range("H63").Value = 3971.5557
$ws.Range("I63").Value = 3348.2856
$ws.Range("K63").Value = 3348.2856
$ws.Range("M63").Value = -2662.2856

$ws.Range("H66").Value = 3971.5557
$ws.Range("I66").Value = 3348.2856
$ws.Range("K66").Value = 16741.428
$ws.Range("M66").Value = -13309.428

$ws.Range("H97").Value = 1193.8889
$ws.Range("I97").Value = 1193.8889
$ws.Range("K97").Value = 1193.8889
$ws.Range("M97").Value = -697.8888999999999

$ws.Range("H102").Value = 22331
$ws.Range("I102").Value = 22331
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 22331
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -20709
$ws.Range("N102").ClearContents()

$ws.Range("H111").Value = 62379.332
$ws.Range("J111").Value = 62379.332
$ws.Range("L111").Value = 62379.332
$ws.Range("N111").Value = -70559.33199999999

$ws.Range("H116").Value = 1012.4
$ws.Range("I116").Value = 1028.2858
$ws.Range("K116").Value = 1028.2858
$ws.Range("M116").Value = 1265.7142

$ws.Range("H132").Value = 5911.8613
$ws.Range("I132").Value = 2505.25
$ws.Range("K132").Value = 7515.75
$ws.Range("M132").Value = -4985.75


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1012.4
$ws.Range("I3").Value = 1028.2858
$ws.Range("K3").Value = 1028.2858
$ws.Range("M3").Value = -914.2858000000001

$ws.Range("H80").Value = 815.5714
$ws.Range("J80").Value = 632.5454999999999
$ws.Range("L80").Value = 632.5454999999999
$ws.Range("N80").Value = -2628.5455

$ws.Range("H83").Value = 815.5714
$ws.Range("J83").Value = 632.5454999999999
$ws.Range("L83").Value = 3162.7275
$ws.Range("N83").Value = -13146.7275

$ws.Range("H98").Value = 79847
$ws.Range("J98").Value = 79847
$ws.Range("L98").Value = 79847
$ws.Range("N98").Value = -85837

$ws.Range("H107").Value = 4416
$ws.Range("I107").Value = 2270
$ws.Range("K107").Value = 2270
$ws.Range("M107").Value = -350


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 126.875
$ws.Range("I2").Value = 98.36364
$ws.Range("J2").Value = 141.80952
$ws.Range("K2").Value = 590.18184
$ws.Range("L2").Value = 850.8571199999999
$ws.Range("M2").Value = -477.18184
$ws.Range("N2").Value = -1076.85712

$ws.Range("H68").Value = 2946.6667
$ws.Range("I68").Value = 2897
$ws.Range("J68").Value = 2971.5
$ws.Range("K68").Value = 8691
$ws.Range("L68").Value = 8914.5
$ws.Range("M68").Value = -7880
$ws.Range("N68").Value = -10536.5

$ws.Range("H71").Value = 2946.6667
$ws.Range("I71").Value = 2897
$ws.Range("J71").Value = 2971.5
$ws.Range("K71").Value = 26073
$ws.Range("L71").Value = 26743.5
$ws.Range("M71").Value = -22017
$ws.Range("N71").Value = -34855.5

$ws.Range("H138").Value = 2684.8
$ws.Range("I138").Value = 2121.1428
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 6363.428400000001
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -1223.428400000001
$ws.Range("N138").Value = -22280

$ws.Range("H140").Value = 117244.16
$ws.Range("I140").Value = 117244.16
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 351732.48
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -346552.48
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 309999.5
$ws.Range("I141").Value = 1004998.3
$ws.Range("J141").Value = 12142.857
$ws.Range("K141").Value = 3014994.9
$ws.Range("L141").Value = 36428.571
$ws.Range("M141").Value = -3009814.9
$ws.Range("N141").Value = -46788.571


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5532
$ws.Range("I122").Value = 4438
$ws.Range("J122").Value = 7720
$ws.Range("K122").Value = 13314
$ws.Range("L122").Value = 23160
$ws.Range("M122").Value = -10864
$ws.Range("N122").Value = -28060


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 29673.334
$ws.Range("I24").Value = 4000
$ws.Range("J24").Value = 42510
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 42510
$ws.Range("M24").Value = -3770
$ws.Range("N24").Value = -42970

$ws.Range("H46").Value = 53331
$ws.Range("J46").Value = 53331
$ws.Range("L46").Value = 53331
$ws.Range("N46").Value = -53793

$ws.Range("H52").Value = 37939.89
$ws.Range("I52").Value = 41819.832
$ws.Range("J52").Value = 30180
$ws.Range("K52").Value = 41819.832
$ws.Range("L52").Value = 30180
$ws.Range("M52").Value = -41593.832
$ws.Range("N52").Value = -30632

$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H113").Value = 883.5217
$ws.Range("I113").Value = 917
$ws.Range("J113").Value = 532
$ws.Range("K113").Value = 2751
$ws.Range("L113").Value = 1596
$ws.Range("M113").Value = -581
$ws.Range("N113").Value = -5936

$ws.Range("H117").Value = 78198
$ws.Range("J117").Value = 78198
$ws.Range("L117").Value = 78198
$ws.Range("N117").Value = -87376

$ws.Range("H124").Value = 110265
$ws.Range("J124").Value = 110265
$ws.Range("L124").Value = 110265
$ws.Range("N124").Value = -120085

$ws.Range("H125").Value = 105974
$ws.Range("J125").Value = 105974
$ws.Range("L125").Value = 105974
$ws.Range("N125").Value = -115814

$ws.Range("H128").Value = 67890
$ws.Range("J128").Value = 67890
$ws.Range("L128").Value = 67890
$ws.Range("N128").Value = -77850

$ws.Range("H134").Value = 53331
$ws.Range("J134").Value = 53331
$ws.Range("L134").Value = 159993
$ws.Range("N134").Value = -165063

$ws.Range("H140").Value = 58464.5
$ws.Range("J140").Value = 58464.5
$ws.Range("L140").Value = 58464.5
$ws.Range("N140").Value = -68824.5
